$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: prepend two new rank-history columns (newest first) before
# the old "B"/"C" columns, pushing the previous two weeks two columns to the
# right (B->D, C->E).
$ws.Range("B1:C1").EntireColumn.Insert()

# Header dates for the two newly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Column C used to carry a fixed width; keep that same width on the (now
# three) date columns C:E it spawned/shifted into.
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 7.1666666666667

# New columns default to "UN" (unchanged) for every analyst row, same as the
# rest of the history columns.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
